# Contrato Promesa de Compraventa 5 - "Cambios al contrato unificado"
#
# Replaces the hard-coded date text ("el quince de marzo del año dos mil
# diecisiete", split across two runs) with the placeholder token
# "plazoVencimiento99", and relocates the document's "_GoBack" bookmark
# (which Word always keeps at the most recently edited spot) from right
# after "NombreRep299" to right after the newly inserted placeholder.

$d = $word.ActiveDocument

# The Find range; after Execute(..., Replace:=wdReplaceAll) it is left
# spanning the text that was just inserted.
$r = $d.Content
$found = $r.Find.Execute(
    "el quince de marzo del año dos mil diecisiete",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "plazoVencimiento99", 2)

if ($found) {
    # Collapse to the end of the freshly-inserted replacement text and
    # move the "_GoBack" bookmark there (Bookmarks.Add on "_GoBack"
    # replaces any existing one, since a document only ever has one).
    $r.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r)
}
